$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 22
$ws.Range("F6").Value = 1152
$ws.Range("F7").Value = 931
$ws.Range("F9").Value = 64
$ws.Range("F10").Value = 82
$ws.Range("F11").Value = 893
$ws.Range("F14").Value = 527
$ws.Range("F16").Value = 123
$ws.Range("F17").Value = 1271
$ws.Range("F18").Value = 2934
$ws.Range("F19").Value = 248
$ws.Range("F20").Value = 1556
$ws.Range("F21").Value = 1306
$ws.Range("F22").Value = 755
$ws.Range("F24").Value = 1307
$ws.Range("F26").Value = 1069
$ws.Range("F28").Value = 3303
$ws.Range("F29").Value = 646
$ws.Range("F30").Value = 550
$ws.Range("F31").Value = 1465

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 6

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 22
$ws.Range("F10").Value = 1152
$ws.Range("F11").Value = 931
$ws.Range("F14").Value = 64
$ws.Range("F20").Value = 82
$ws.Range("F21").Value = 6
$ws.Range("F23").Value = 893
$ws.Range("F26").Value = 527
$ws.Range("F28").Value = 123
$ws.Range("F29").Value = 1271
$ws.Range("F30").Value = 2934
$ws.Range("F31").Value = 248
$ws.Range("F32").Value = 1556
$ws.Range("F33").Value = 1306
$ws.Range("F34").Value = 755
$ws.Range("F36").Value = 1307
$ws.Range("F40").Value = 1069
$ws.Range("F42").Value = 3303
$ws.Range("F43").Value = 646
$ws.Range("F44").Value = 550
$ws.Range("F45").Value = 1465
